$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.410.98"
$ws.Range("E2").Value = "  -0.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.080.64"
$ws.Range("E3").Value = "  +5.41%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "199.21"
$ws.Range("E5").Value = "  +0.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "618.13"
$ws.Range("E6").Value = "  +4.11%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("E8").Value = "  +0.90%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.213"
$ws.Range("E9").Value = "  +7.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.077.30"
$ws.Range("E10").Value = "  +5.44%  "

$ws.Range("E11").Value = "  +2.24%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.25"
$ws.Range("E13").Value = "  +8.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.644.89"
$ws.Range("E14").Value = "  +5.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.38"
$ws.Range("E15").Value = "  +4.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000196"
$ws.Range("E16").Value = "  +4.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "76.201.69"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.079.09"
$ws.Range("E18").Value = "  +5.30%  "

$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.02"
$ws.Range("E20").Value = "  +3.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.59"
$ws.Range("E21").Value = "  +15.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.05"
$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.51"
$ws.Range("E23").Value = "  +5.61%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.45"
$ws.Range("E24").Value = "  +0.86%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.63"
$ws.Range("E25").Value = "  +9.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.232.46"
$ws.Range("E26").Value = "  +5.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.48"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.08"
$ws.Range("E29").Value = "  +5.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000109"
$ws.Range("E30").Value = "  +2.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.996"
$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.33"
$ws.Range("E32").Value = "  +0.78%  "

$ws.Range("E33").Value = "  +4.84%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "502.82"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("E35").Value = "  +6.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.126"
$ws.Range("E37").Value = "  +15.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.89"
$ws.Range("E38").Value = "  +4.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "163.42"
$ws.Range("E39").Value = "  -1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "195.44"
$ws.Range("E40").Value = "  +9.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.07"
$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.379"
$ws.Range("E42").Value = "  -2.68%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.103"
$ws.Range("E43").Value = "  -6.93%  "

$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.805"
$ws.Range("E45").Value = "  +23.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.18"
$ws.Range("E46").Value = "  +5.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.27"
$ws.Range("E47").Value = "  +7.51%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.67"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("E49").Value = "  +6.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.73"
$ws.Range("E50").Value = "  +1.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.599"
$ws.Range("E51").Value = "  +2.02%  "
